$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6.. down by one.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly data point.
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = "2022-05-10"
$ws.Cells.Item(6, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112035
$ws.Cells.Item(6, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 55
$ws.Cells.Item(6, 11).Value = 30000
$ws.Cells.Item(6, 12).Value = 30000
$ws.Cells.Item(6, 13).Value = 30000
$ws.Cells.Item(6, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(6, 16).Value = 3000
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = "Hortaliza"
